$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12:B12").Copy($ws.Range("A13:B13"))

$ws.Range("A13").Value = "ItemData"
$ws.Range("B13").Value = "ItemData.xlsx"
